# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped values. Source cells are stored as text (e.g.
# "605.68", "  +0.29%  "), so for any new value that looks like a plain
# number we briefly force a text number format before assigning it and
# clear the format again afterwards - otherwise Range.Value would silently
# coerce strings such as "1.00" / "2.50" into the numbers 1 / 2.5 (losing
# the trailing zero) or introduce floating-point noise (e.g. 605.03 ->
# 605.02999999999997).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.355.82"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.553.39"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "3.551.22"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.87"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "4.158.42"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.02"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.542.92"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "66.462.76"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.84"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.07"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.54"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "3.695.23"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.50"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.95"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.42"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "3.549.56"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.154"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.65"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0847"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.03"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +5.28%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.05"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.83%  "
